$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 8; existing rows 8.. shift down to 10..,
# keeping formatting (date style on column D) intact.
$ws.Rows("8:9").Insert()

# New row 8: Especial, 2022-04-19
$ws.Cells.Item(8, 1).Value = 4
$ws.Cells.Item(8, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(8, 3).Value = "Los Lagos"
$ws.Cells.Item(8, 4).Value = (Get-Date -Year 2022 -Month 4 -Day 19).Date
$ws.Cells.Item(8, 5).Value = 10
$ws.Cells.Item(8, 6).Value = "Fruta"
$ws.Cells.Item(8, 7).Value = 100104
$ws.Cells.Item(8, 8).Value = "Frutos de pepita"
$ws.Cells.Item(8, 9).Value = 100104003
$ws.Cells.Item(8, 10).Value = "Membrillo"
$ws.Cells.Item(8, 11).Value = "Champion"
$ws.Cells.Item(8, 12).Value = "Especial"
$ws.Cells.Item(8, 13).Value = 150
$ws.Cells.Item(8, 14).Value = 18000
$ws.Cells.Item(8, 15).Value = 18000
$ws.Cells.Item(8, 16).Value = 18000
$ws.Cells.Item(8, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(8, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(8, 19).Value = 1000
$ws.Cells.Item(8, 20).Value = 18

# New row 9: Primera, 2022-04-19
$ws.Cells.Item(9, 1).Value = 4
$ws.Cells.Item(9, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(9, 3).Value = "Los Lagos"
$ws.Cells.Item(9, 4).Value = (Get-Date -Year 2022 -Month 4 -Day 19).Date
$ws.Cells.Item(9, 5).Value = 10
$ws.Cells.Item(9, 6).Value = "Fruta"
$ws.Cells.Item(9, 7).Value = 100104
$ws.Cells.Item(9, 8).Value = "Frutos de pepita"
$ws.Cells.Item(9, 9).Value = 100104003
$ws.Cells.Item(9, 10).Value = "Membrillo"
$ws.Cells.Item(9, 11).Value = "Champion"
$ws.Cells.Item(9, 12).Value = "Primera"
$ws.Cells.Item(9, 13).Value = 150
$ws.Cells.Item(9, 14).Value = 16000
$ws.Cells.Item(9, 15).Value = 16000
$ws.Cells.Item(9, 16).Value = 16000
$ws.Cells.Item(9, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(9, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(9, 19).Value = 889
$ws.Cells.Item(9, 20).Value = 18
